# "update data with resort sheetname"
#
# The workbook currently lists the sheets as:
#   1) "2022-Q2"  (per-fund holding detail table, A1:H27)
#   2) "总计"      (quarter summary totals, A1:D2)
#
# The sheets need to be re-sorted so that the summary sheet "总计" comes
# first and the quarterly detail sheet "2022-Q2" comes second - i.e. the
# two tabs simply swap places (their cell contents/styles are untouched).

$wb = $excel.ActiveWorkbook

$ws2022 = $wb.Worksheets.Item("2022-Q2")
$wsTotal = $wb.Worksheets.Item("总计")

# Move "总计" so that it sits right before "2022-Q2", i.e. it becomes the
# first tab and "2022-Q2" becomes the second tab.
$wsTotal.Move($ws2022)

# "2022-Q2" was the selected/active tab before the re-sort; keep it that
# way after reordering (look it up again by name since the worksheet
# collection positions have changed).
$wb.Worksheets.Item("2022-Q2").Activate()
